# feat: add 2022-Q3 data
#
# Before: sheet1 "总计" (totals) + sheet2 "2022-Q2" (fund holdings for Q2)
# After : sheet1 "总计" gets a new leading row for 2022-Q3 (old Q2 row shifts
#         down); the old "2022-Q2" worksheet is duplicated - the duplicate
#         keeps the "2022-Q2" name/data, while the original sheet is
#         renamed "2022-Q3" and its data is replaced with the new Q3 fund
#         holdings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for 2022-Q3, push 2022-Q2 to row 3
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

$totals.Rows.Item(2).Insert()

# Row 3 (formerly row 2, "2022-Q2") keeps its data/format; just fix its
# running index (A column) back to 1.
$totals.Range("A3").Value = 1

# New row 2 ("2022-Q3") - clone the style of A3 onto A2 so it matches the
# sheet's numbered-row styling, then fill in the values.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 0.07

# ---------------------------------------------------------------------
# 2. Duplicate the "2022-Q2" worksheet so the old fund-holdings data is
#    preserved under its original name, directly after itself.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($null, $q2)

$q2Copy = $wb.Worksheets.Item(3)
$q2Copy.Name = "2022-Q2-new"
$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 3. Replace the (now renamed) "2022-Q3" sheet's data with the new
#    quarter's fund-holdings table.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)

# Drop the now-unneeded trailing rows (old sheet had 10 data rows, new
# one only has 5) so the used range shrinks back down.
$q3.Range("A7:H11").Delete()

# Match header + A-column styling to the "s=2" style used on the totals
# sheet's header/number cells.
$totals.Range("B1:D1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$totals.Range("B1").Copy()
$q3.Range("A2:A6").PasteSpecial(-4122)

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'008602"
$q3.Range("C2").Value = "方正富邦新兴成长混合A"
$q3.Range("D2").Value = "'1.23"
$q3.Range("E2").Value = "'86.03"
$q3.Range("F2").Value = "'4.09"
$q3.Range("G2").Value = "'0.0503"
$q3.Range("H2").Value = 5

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'004332"
$q3.Range("C3").Value = "恒生前海沪港深新兴产业精选混合"
$q3.Range("D3").Value = "'0.47"
$q3.Range("E3").Value = "'92.74"
$q3.Range("F3").Value = "'2.84"
$q3.Range("G3").Value = "'0.0133"
$q3.Range("H3").Value = 10

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'006347"
$q3.Range("C4").Value = "安信量化优选股票C"
$q3.Range("D4").Value = "'0.15"
$q3.Range("E4").Value = "'90.50"
$q3.Range("F4").Value = "'0.84"
$q3.Range("G4").Value = "'0.0013"
$q3.Range("H4").Value = 4

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'008603"
$q3.Range("C5").Value = "方正富邦新兴成长混合C"
$q3.Range("D5").Value = "'0.03"
$q3.Range("E5").Value = "'86.03"
$q3.Range("F5").Value = "'4.09"
$q3.Range("G5").Value = "'0.0012"
$q3.Range("H5").Value = 5

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "'006346"
$q3.Range("C6").Value = "安信量化优选股票A"
$q3.Range("D6").Value = "'0.03"
$q3.Range("E6").Value = "'90.50"
$q3.Range("F6").Value = "'0.84"
$q3.Range("G6").Value = "'0.0003"
$q3.Range("H6").Value = 4
